# Auto-generated COM-interop edit script
# Adds rows 42-61 to the 'Orders' sheet and updates Summary!G2

$wb = $excel.ActiveWorkbook
$orders = $wb.Worksheets.Item("Orders")
$summary = $wb.Worksheets.Item("Summary")

# New order rows (PackageID in col A, FlowerName in col C, Number in col F).
# NumberFormat is forced to Text ("@") before assigning any digit-only string so
# Excel keeps it stored as text (matching the rest of the sheet) instead of coercing
# it to a numeric value.
$newOrderRows = @(
    @{ Row = 42; A = '8'; C = '589_洋牡丹香槟_undefined_undefined_1bunch'; F = '10' }
    @{ Row = 43; A = $null; C = '586_洋牡丹白_undefined_undefined_1bunch'; F = '20' }
    @{ Row = 44; A = $null; C = '590_洋牡丹粉_undefined_undefined_1bunch'; F = '15' }
    @{ Row = 45; A = $null; C = '587_洋牡丹橙_undefined_undefined_1bunch'; F = '5' }
    @{ Row = 46; A = $null; C = '588_洋牡丹黄_undefined_undefined_1bunch'; F = '5' }
    @{ Row = 47; A = $null; C = '591_洋牡丹黑_undefined_undefined_1bunch'; F = '5' }
    @{ Row = 48; A = $null; C = '648_洋牡丹河内_undefined_undefined_1bunch'; F = '10' }
    @{ Row = 49; A = $null; C = '649_洋牡丹樱花粉_undefined_undefined_1bunch'; F = '10' }
    @{ Row = 50; A = '9'; C = '651_大丽花 奶油桃子_undefined_undefined_5stems'; F = '15' }
    @{ Row = 51; A = $null; C = '653_大丽花 黑_undefined_undefined_5stems'; F = '5' }
    @{ Row = 52; A = $null; C = '656_大丽花 梅根_undefined_undefined_5stems'; F = '10' }
    @{ Row = 53; A = $null; C = '551_铁线莲_Glematis_undefined_1bunch'; F = '5' }
    @{ Row = 54; A = $null; C = '413_风铃花淡紫色_Canterbury Bells
light purple_undefined_1bunch'; F = '5' }
    @{ Row = 55; A = '10'; C = '144_高原红_High Plateau Red_Rosa rugosa Thunb._20stems'; F = '15' }
    @{ Row = 56; A = $null; C = '41_拉丝白_Spider White_Gerbera L._20stems'; F = '10' }
    @{ Row = 57; A = $null; C = '46_拉丝橙_Spider orange_Gerbera L._20stems'; F = '5' }
    @{ Row = 58; A = $null; C = '630_吸色康乃馨天蓝_tinted tiffany blue_undefined_20stems'; F = '10' }
    @{ Row = 59; A = $null; C = '631_吸色康乃馨宝蓝_tinted blue_undefined_20stems'; F = '10' }
    @{ Row = 60; A = $null; C = '508_风铃花白色_Canterbury Bells 
white_undefined_1bunch'; F = '10' }
    @{ Row = 61; A = $null; C = '414_风铃花粉色_Canterbury Bells
pink_undefined_1bunch'; F = $null }
)

foreach ($r in $newOrderRows) {
    if ($null -ne $r.A) {
        $cell = $orders.Cells.Item($r.Row, 1)
        $cell.NumberFormat = "@"
        $cell.Value = $r.A
    }
    if ($null -ne $r.C) {
        $orders.Cells.Item($r.Row, 3).Value = $r.C
    }
    if ($null -ne $r.F) {
        $cell = $orders.Cells.Item($r.Row, 6)
        $cell.NumberFormat = "@"
        $cell.Value = $r.F
    }
}

# Update the rolled-up Number column string on the Summary sheet (G2) to include
# the Number values of the newly added order rows. Force text storage so the long
# digit string isn't rounded into scientific notation.
$g2 = $summary.Cells.Item(2, 7)
$g2.NumberFormat = "@"
$g2.Value = '0151540401033532151014713101491410105510115111082615151515151041595010201555510101551055151051010100'

